$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Header row: new columns J (…_A) and K (…_0) for the 04-10 snapshot ----
$ws.Cells.Item(1, 10).Value2 = "04-10_A"
$ws.Cells.Item(1, 11).Value2 = "04-10_0"
$ws.Range("J1:K1").Font.Bold = $true
$ws.Range("J1:K1").Borders.LineStyle = 1
$ws.Range("J1:K1").HorizontalAlignment = -4108
$ws.Range("J1:K1").VerticalAlignment = -4160

# ---- Per-row data ----
# color: Interior.Color to copy from column H onto the new column J (null = no fill)
# hval : column H's numeric value, copied verbatim into column J
# kind : "text"     -> column I held this week's text record; convert it to a real
#                       number and move the text record into the new column K
#        "blank"    -> column I was already numeric (no new record this week);
#                       the new J/K columns stay blank
#        "emptyrow" -> the whole row had no data at all (blank placeholder row)
# ival : the inline-text value that used to live in column I (becomes column K)
$rows = @(
    @{ r=2; color=255; hval=0; kind="blank"; ival="0.0" },
    @{ r=3; color=255; hval=0; kind="text"; ival="2519" },
    @{ r=4; color=255; hval=0; kind="text"; ival="0" },
    @{ r=5; color=255; hval=0; kind="text"; ival="2498" },
    @{ r=6; color=255; hval=0; kind="text"; ival="0" },
    @{ r=7; color=255; hval=0; kind="text"; ival="2500" },
    @{ r=8; color=32768; hval=40; kind="blank"; ival="4975.0" },
    @{ r=9; color=65535; hval=8; kind="text"; ival="2625" },
    @{ r=10; color=255; hval=0; kind="text"; ival="0" },
    @{ r=11; color=255; hval=0; kind="text"; ival="0" },
    @{ r=12; color=255; hval=0; kind="text"; ival="0" },
    @{ r=13; color=255; hval=0; kind="text"; ival="2491" },
    @{ r=14; color=255; hval=0; kind="text"; ival="2499" },
    @{ r=15; color=16777215; hval=29; kind="text"; ival="2778" },
    @{ r=16; color=255; hval=0; kind="text"; ival="2500" },
    @{ r=17; color=32768; hval=40; kind="text"; ival="2956" },
    @{ r=18; color=16777215; hval=20; kind="text"; ival="2740" },
    @{ r=19; color=16777215; hval=23; kind="text"; ival="2904" },
    @{ r=20; color=32768; hval=33; kind="text"; ival="3040" },
    @{ r=21; color=16777215; hval=30; kind="text"; ival="2929" },
    @{ r=22; color=255; hval=0; kind="text"; ival="2498" },
    @{ r=23; color=16777215; hval=20; kind="text"; ival="2741" },
    @{ r=24; color=255; hval=0; kind="text"; ival="0" },
    @{ r=25; color=255; hval=0; kind="text"; ival="2526" },
    @{ r=26; color=32768; hval=34; kind="text"; ival="2967" },
    @{ r=27; color=65535; hval=10; kind="text"; ival="2672" },
    @{ r=28; color=255; hval=0; kind="text"; ival="2500" },
    @{ r=29; color=16777215; hval=20; kind="text"; ival="2830" },
    @{ r=30; color=255; hval=0; kind="text"; ival="2524" },
    @{ r=31; color=16777215; hval=30; kind="text"; ival="2906" },
    @{ r=32; color=16777215; hval=30; kind="text"; ival="2856" },
    @{ r=33; color=16777215; hval=30; kind="text"; ival="3052" },
    @{ r=34; color=16777215; hval=21; kind="text"; ival="2738" },
    @{ r=35; color=65535; hval=14; kind="text"; ival="2767" },
    @{ r=36; color=255; hval=0; kind="text"; ival="2590" },
    @{ r=37; color=16777215; hval=30; kind="text"; ival="2904" },
    @{ r=38; color=65535; hval=5; kind="text"; ival="2528" },
    @{ r=39; color=65535; hval=5; kind="text"; ival="2575" },
    @{ r=40; color=16777215; hval=30; kind="text"; ival="2991" },
    @{ r=41; color=16777215; hval=30; kind="text"; ival="2997" },
    @{ r=42; color=16777215; hval=30; kind="text"; ival="2990" },
    @{ r=43; color=16777215; hval=23; kind="text"; ival="3024" },
    @{ r=44; color=65535; hval=13; kind="text"; ival="2710" },
    @{ r=45; color=16777215; hval=30; kind="text"; ival="2933" },
    @{ r=46; color=65535; hval=12; kind="text"; ival="2706" },
    @{ r=47; color=255; hval=0; kind="text"; ival="2567" },
    @{ r=48; color=16777215; hval=20; kind="text"; ival="2915" },
    @{ r=49; color=16777215; hval=20; kind="text"; ival="2792" },
    @{ r=50; color=65535; hval=5; kind="blank"; ival="4139.0" },
    @{ r=51; color=16777215; hval=23; kind="text"; ival="2791" },
    @{ r=52; color=16777215; hval=20; kind="text"; ival="2786" },
    @{ r=53; color=16777215; hval=20; kind="text"; ival="2941" },
    @{ r=54; color=16777215; hval=20; kind="text"; ival="2745" },
    @{ r=55; color=16777215; hval=30; kind="text"; ival="2845" },
    @{ r=56; color=16777215; hval=21; kind="text"; ival="2794" },
    @{ r=57; color=16777215; hval=20; kind="text"; ival="2924" },
    @{ r=58; color=255; hval=0; kind="text"; ival="0" },
    @{ r=59; color=255; hval=0; kind="text"; ival="0" },
    @{ r=60; color=255; hval=0; kind="text"; ival="2500" },
    @{ r=61; color=255; hval=0; kind="text"; ival="2513" },
    @{ r=62; color=255; hval=0; kind="text"; ival="0" },
    @{ r=63; color=16777215; hval=$null; kind="emptyrow"; ival=$null },
    @{ r=64; color=255; hval=0; kind="text"; ival="0" },
    @{ r=65; color=16777215; hval=20; kind="text"; ival="2796" },
    @{ r=66; color=255; hval=0; kind="text"; ival="2498" },
    @{ r=67; color=32768; hval=38; kind="text"; ival="3051" },
    @{ r=68; color=255; hval=0; kind="text"; ival="0" },
    @{ r=69; color=255; hval=0; kind="text"; ival="0" },
    @{ r=70; color=255; hval=0; kind="text"; ival="0" },
    @{ r=71; color=255; hval=0; kind="text"; ival="2539" },
    @{ r=72; color=255; hval=0; kind="text"; ival="2505" },
    @{ r=73; color=255; hval=0; kind="text"; ival="0" },
    @{ r=74; color=255; hval=0; kind="text"; ival="0" },
    @{ r=75; color=255; hval=0; kind="text"; ival="0" },
    @{ r=76; color=255; hval=0; kind="text"; ival="1500" },
    @{ r=77; color=65535; hval=4; kind="text"; ival="2574" },
    @{ r=78; color=255; hval=0; kind="text"; ival="0" },
    @{ r=79; color=255; hval=0; kind="text"; ival="0" },
    @{ r=80; color=255; hval=0; kind="text"; ival="0" },
    @{ r=81; color=255; hval=0; kind="text"; ival="0" },
    @{ r=82; color=255; hval=0; kind="text"; ival="0" },
    @{ r=83; color=255; hval=0; kind="text"; ival="0" },
    @{ r=84; color=255; hval=0; kind="text"; ival="0" },
    @{ r=85; color=255; hval=0; kind="text"; ival="0" },
    @{ r=86; color=255; hval=0; kind="text"; ival="0" },
    @{ r=87; color=255; hval=0; kind="text"; ival="0" },
    @{ r=88; color=255; hval=0; kind="text"; ival="2499" },
    @{ r=89; color=255; hval=0; kind="text"; ival="0" },
    @{ r=90; color=255; hval=0; kind="text"; ival="0" },
    @{ r=91; color=255; hval=0; kind="text"; ival="2499" },
    @{ r=92; color=255; hval=0; kind="text"; ival="0" },
    @{ r=93; color=255; hval=0; kind="text"; ival="0" },
    @{ r=94; color=255; hval=0; kind="text"; ival="0" },
    @{ r=95; color=255; hval=0; kind="text"; ival="0" },
    @{ r=96; color=255; hval=0; kind="text"; ival="0" },
    @{ r=97; color=255; hval=0; kind="text"; ival="0" },
    @{ r=98; color=255; hval=0; kind="text"; ival="0" },
    @{ r=99; color=255; hval=0; kind="text"; ival="0" },
    @{ r=100; color=255; hval=0; kind="text"; ival="0" },
    @{ r=101; color=255; hval=0; kind="text"; ival="1499" },
    @{ r=102; color=255; hval=0; kind="text"; ival="0" },
    @{ r=103; color=255; hval=0; kind="text"; ival="0" },
    @{ r=104; color=255; hval=0; kind="text"; ival="0" },
    @{ r=105; color=255; hval=0; kind="text"; ival="0" },
    @{ r=106; color=255; hval=0; kind="text"; ival="0" },
    @{ r=107; color=255; hval=0; kind="text"; ival="0" },
    @{ r=108; color=255; hval=0; kind="text"; ival="0" },
    @{ r=109; color=255; hval=0; kind="text"; ival="0" },
    @{ r=110; color=255; hval=0; kind="text"; ival="0" },
    @{ r=111; color=255; hval=0; kind="text"; ival="0" },
    @{ r=112; color=16777215; hval=20; kind="text"; ival="2831" },
    @{ r=113; color=255; hval=0; kind="text"; ival="2519" }
)

foreach ($row in $rows) {
    $r = $row.r

    if ($row.color -ne $null) {
        $ws.Cells.Item($r, 10).Interior.Color = $row.color
    }

    if ($row.kind -eq "text") {
        $ws.Cells.Item($r, 9).Value2 = [double]$row.ival
        if ($row.hval -ne $null) {
            $ws.Cells.Item($r, 10).Value2 = $row.hval
        }
        $ws.Cells.Item($r, 11).Value = "'" + $row.ival
    } elseif ($row.kind -eq "blank") {
        if ($row.hval -ne $null) {
            $ws.Cells.Item($r, 10).Value2 = $row.hval
        }
    }
}

# ---- Rows 112/113: the trailing placeholder rows store their id as a real number now ----
$ws.Cells.Item(112, 1).Value2 = 27484940
$ws.Cells.Item(113, 1).Value2 = 41837764
